# Weekly dashboard data refresh - 2026-01-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RGDP (row 3-4) ---
$ws.Range("F3").Value = 24026.834
$ws.Range("F4").Value = 0.0107634621312982

# --- RPCE (row 9-10) ---
$ws.Range("F9").Value = 16585.878
$ws.Range("F10").Value = 0.008524606910566446

# --- UI Initial Claims (row 13): new week added, prior weeks shift right ---
$ws.Range("N13").Value = 46034
$ws.Range("U13").Value = 215000
$ws.Range("T13").Value = 200000
$ws.Range("S13").Value = 207000
$ws.Range("R13").Value = 199000
$ws.Range("Q13").Value = 200000

# --- UI Continuing Claims (row 14): new week added, prior weeks shift right ---
$ws.Range("N14").Value = 46027
$ws.Range("U14").Value = 1914000
$ws.Range("T14").Value = 1856000
$ws.Range("S14").Value = 1903000
$ws.Range("R14").Value = 1875000
$ws.Range("Q14").Value = 1849000

# --- Gross Priv Fixed Inv Non Res / Res (row 26-27) ---
$ws.Range("F26").Value = 0.007794617579222285
$ws.Range("F27").Value = -0.01832684674964358

# --- 5yr,5yr Forward (row 29): new day added, prior days shift right ---
$ws.Range("N29").Value = 46043
$ws.Range("U29").Value = 2.22
$ws.Range("T29").Value = 2.22
$ws.Range("S29").Value = 2.27
$ws.Range("R29").Value = 2.26

# --- 10yr TIPS (row 30): new day added, prior days shift right ---
$ws.Range("N30").Value = 46043
$ws.Range("U30").Value = 2.29
$ws.Range("T30").Value = 2.29
$ws.Range("S30").Value = 2.33
$ws.Range("R30").Value = 2.33
$ws.Range("Q30").Value = 2.34

# --- Gov. Cons (row 44-45) ---
$ws.Range("F44").Value = 5324.402
$ws.Range("F45").Value = 0.01669515005814426

# --- FFR (row 47) ---
$ws.Range("N47").Value = 46042

# --- 2y UST (row 48): new day added, prior days shift right ---
$ws.Range("N48").Value = 46042
$ws.Range("U48").Value = 3.53
$ws.Range("T48").Value = 3.51
$ws.Range("S48").Value = 3.56
$ws.Range("R48").Value = 3.59
$ws.Range("Q48").Value = 3.6

# --- 5y UST (row 49): new day added, prior days shift right ---
$ws.Range("N49").Value = 46042
$ws.Range("U49").Value = 3.75
$ws.Range("T49").Value = 3.72
$ws.Range("S49").Value = 3.77
$ws.Range("R49").Value = 3.82
$ws.Range("Q49").Value = 3.86

# --- 10y UST (row 50): new day added, prior days shift right ---
$ws.Range("N50").Value = 46042
$ws.Range("U50").Value = 4.18
$ws.Range("T50").Value = 4.15
$ws.Range("S50").Value = 4.17
$ws.Range("R50").Value = 4.24
$ws.Range("Q50").Value = 4.3

# --- BAA (row 52): new day added, prior days shift right ---
$ws.Range("N52").Value = 46042
$ws.Range("U52").Value = 5.87
$ws.Range("T52").Value = 5.83
$ws.Range("S52").Value = 5.82
$ws.Range("R52").Value = 5.87
$ws.Range("Q52").Value = 5.95

# --- Drop stale "new data" yellow highlight on N41,N42,N43,N44,N51 ---
# (copy the unhighlighted format from N3, which already carries that style)
$ws.Range("N3").Copy()
$ws.Range("N41").PasteSpecial(-4122)
$ws.Range("N42").PasteSpecial(-4122)
$ws.Range("N43").PasteSpecial(-4122)
$ws.Range("N44").PasteSpecial(-4122)
$ws.Range("N51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
